# Germany Landesliga - atualizacao de bases, 30-03-2024 19:32
#
# The underlying source data re-sorted a couple of match rows; the team
# names / ids / scores / odds for those matches moved between rows while
# the row's id / Div / Div Original Name / Date (columns A, C, D, E)
# stayed put. Reproduce that by swapping/rotating the B..AC payload of the
# affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2 and 3: swap the match payload (columns B, F..AC) -------------
$row2 = $ws.Range("B2:AC2")
$row3 = $ws.Range("B3:AC3")

$v2 = $row2.Value()
$v3 = $row3.Value()

$row2.Value = $v3
$row3.Value = $v2

# --- Rows 11, 12, 13: rotate the match payload -----------------------
# row11 <- row12, row12 <- row13, row13 <- row11
$row11 = $ws.Range("B11:AC11")
$row12 = $ws.Range("B12:AC12")
$row13 = $ws.Range("B13:AC13")

$v11 = $row11.Value()
$v12 = $row12.Value()
$v13 = $row13.Value()

$row11.Value = $v12
$row12.Value = $v13
$row13.Value = $v11
